# PAS-730: Small fixed for old test faced during execution.
# Changed dates for control tables to delete overlap.
#
# The two "control table" rows (row 2 / row 3) used placeholder text codes
# "039" and "040" (row 3, columns AA/AB) and a leftover "SS" test symbol
# (row 3, columns AE:AH). Replace them with the corrected values so the
# unused shared strings ("039", "040") disappear and "SS" becomes "X".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AA3 / AB3: was text "039" / "040" -> now plain numbers 39 / 40
$ws.Range("AA3").Value = 39
$ws.Range("AB3").Value = 40

# AE3:AH3: was "SS" -> now "X"
$ws.Range("AE3").Value = "X"
$ws.Range("AF3").Value = "X"
$ws.Range("AG3").Value = "X"
$ws.Range("AH3").Value = "X"

# Restore selection to AH3 (also clears the stale topLeftCell="Z1" scroll position)
[void]$ws.Range("AH3").Select()
